$wb = $excel.ActiveWorkbook

# --- Suite1: Cart page run flag Y -> N ---
$wsSuite1 = $wb.Worksheets.Item("Suite1")
$wsSuite1.Range("B3").Value = "N"

# --- Suite4: Order Detail Page row6 run flag N -> Y ---
$wsSuite4 = $wb.Worksheets.Item("Suite4")
$wsSuite4.Range("B6").Value = "Y"

# --- AppControl: Order Detail Page (Suite4) run flag N -> Y ---
$wsApp = $wb.Worksheets.Item("AppControl")
$wsApp.Range("B6").Value = "Y"

# --- Update the selections (active cell) on each affected sheet ---
# (AppControl is selected last so it remains the active tab, matching
#  the original workbook where AppControl's sheetView has tabSelected="1")
$wsSuite1.Range("C6").Select() | Out-Null
$wsSuite4.Range("C9").Select() | Out-Null
$wsApp.Range("D5").Select() | Out-Null
